$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Prepare the destination formatting by copying the existing "data row" /
# "running summary row" pair (rows 47-48) down onto the two new pairs
# (50-51 and 53-54), and copying the trailing "C52/D52" summary cell's
# formatting onto its new home at C57/D57 - all BEFORE we touch row 52 so
# the source formatting is still intact when we copy it.
# ---------------------------------------------------------------------------
$ws.Range("C52:D52").Copy()
$ws.Range("C57:D57").PasteSpecial(-4122)

$ws.Range("A47:J48").Copy()
$ws.Range("A50:J51").PasteSpecial(-4122)
$ws.Range("A53:J54").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 50: Repair to Wing Attach (new data row)
# ---------------------------------------------------------------------------
$ws.Range("A50").Value = 42225
$ws.Range("B50").Value = "Repair to Wing Attach"
$ws.Range("C50").Value = 40
$ws.Range("D50").Value = -30
$ws.Range("E50").Formula = "=E35"
$ws.Range("F50").Formula = "=F35"
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = "Mass estimated based on position estimate and ballast required to balance"

# ---------------------------------------------------------------------------
# Row 51: running summary after row 50
# ---------------------------------------------------------------------------
$ws.Range("B51").Value = "Airframe"
$ws.Range("C51").Formula = "=C48+C50"
$ws.Range("D51").Formula = "=(D48*`$C48 + D50*`$C50)/`$C51"
$ws.Range("E51").Formula = "=(E48*`$C48 + E50*`$C50)/`$C51"
$ws.Range("F51").Formula = "=(F48*`$C48 + F50*`$C50)/`$C51"
$ws.Range("G51").Formula = "=(G48+(`$C48/1000)*((`$E48*in2m-`$E51*in2m)^2+(`$F48*in2m-`$F51*in2m)^2)) + SIGN(`$C50)*((G50)+ABS(`$C50/1000)*((`$E50*in2m-`$E51*in2m)^2+(`$F50*in2m-`$F51*in2m)^2))"
$ws.Range("H51").Formula = "=(H48+(`$C48/1000)*((`$D48*in2m-`$D51*in2m)^2+(`$F48*in2m-`$F51*in2m)^2)) + SIGN(`$C50)*((H50)+ABS(`$C50/1000)*((`$D50*in2m-`$D51*in2m)^2+(`$F50*in2m-`$F51*in2m)^2))"
$ws.Range("I51").Formula = "=(I48+(`$C48/1000)*((`$D48*in2m-`$D51*in2m)^2+(`$E48*in2m-`$E51*in2m)^2)) + SIGN(`$C50)*((I50)+ABS(`$C50/1000)*((`$D50*in2m-`$D51*in2m)^2+(`$E50*in2m-`$E51*in2m)^2))"

# ---------------------------------------------------------------------------
# Row 53: Added Ballast (new data row)
# ---------------------------------------------------------------------------
$ws.Range("A53").Value = 42225
$ws.Range("B53").Value = "Added Ballast"
$ws.Range("C53").Value = 15
$ws.Range("D53").Value = -4.62
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = -0.9
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = "Assume negligible body inertia"

# ---------------------------------------------------------------------------
# Row 54: running summary after row 53
# ---------------------------------------------------------------------------
$ws.Range("B54").Value = "Airframe"
$ws.Range("C54").Formula = "=C51+C53"
$ws.Range("D54").Formula = "=(D51*`$C51 + D53*`$C53)/`$C54"
$ws.Range("E54").Formula = "=(E51*`$C51 + E53*`$C53)/`$C54"
$ws.Range("F54").Formula = "=(F51*`$C51 + F53*`$C53)/`$C54"
$ws.Range("G54").Formula = "=(G51+(`$C51/1000)*((`$E51*in2m-`$E54*in2m)^2+(`$F51*in2m-`$F54*in2m)^2)) + SIGN(`$C53)*((G53)+ABS(`$C53/1000)*((`$E53*in2m-`$E54*in2m)^2+(`$F53*in2m-`$F54*in2m)^2))"
$ws.Range("H54").Formula = "=(H51+(`$C51/1000)*((`$D51*in2m-`$D54*in2m)^2+(`$F51*in2m-`$F54*in2m)^2)) + SIGN(`$C53)*((H53)+ABS(`$C53/1000)*((`$D53*in2m-`$D54*in2m)^2+(`$F53*in2m-`$F54*in2m)^2))"
$ws.Range("I54").Formula = "=(I51+(`$C51/1000)*((`$D51*in2m-`$D54*in2m)^2+(`$E51*in2m-`$E54*in2m)^2)) + SIGN(`$C53)*((I53)+ABS(`$C53/1000)*((`$D53*in2m-`$D54*in2m)^2+(`$E53*in2m-`$E54*in2m)^2))"
$ws.Range("J54").Value = "Balance performed to establish required ballast mass."

# ---------------------------------------------------------------------------
# Relocate the trailing lb-conversion summary from row 52 to row 57
# ---------------------------------------------------------------------------
$ws.Range("C57").Formula = "=C48/1000*2.20462"
$ws.Range("D57").Value = "lb"

# Remove the old row 52 entirely (contents + formatting) now that its
# content lives on at row 57.
$ws.Range("C52:D52").Clear()

